$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt9a"
$ws.Cells.Item(2, 3).Value = "Fzd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.6836936666666666
$ws.Cells.Item(2, 8).Value = 2.051081
$ws.Cells.Item(2, 9).Value = 0.07240868516880868
$ws.Cells.Item(2, 10).Value = 0.07240868516880865
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 13.57958433333333
$ws.Cells.Item(2, 14).Value = 40.738753
$ws.Cells.Item(2, 15).Value = 0.2289698008477291
$ws.Cells.Item(2, 16).Value = 0.2289698008477291
$ws.Cells.Item(2, 17).Value = 9.28427580466589
$ws.Cells.Item(2, 18).Value = 83.558482241993
$ws.Cells.Item(2, 19).Value = 0.01657940222274804
$ws.Cells.Item(2, 20).Value = 0.01657940222274803

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt9a"
$ws.Cells.Item(3, 3).Value = "Fzd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.6836936666666666
$ws.Cells.Item(3, 8).Value = 2.051081
$ws.Cells.Item(3, 9).Value = 0.07240868516880868
$ws.Cells.Item(3, 10).Value = 0.07240868516880865
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 19.768727
$ws.Cells.Item(3, 14).Value = 59.306181
$ws.Cells.Item(3, 15).Value = 0.3333269541315948
$ws.Cells.Item(3, 16).Value = 0.3333269541315948
$ws.Cells.Item(3, 17).Value = 13.51575344796233
$ws.Cells.Item(3, 18).Value = 121.641781031661
$ws.Cells.Item(3, 19).Value = 0.02413576647999258
$ws.Cells.Item(3, 20).Value = 0.02413576647999257

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt9a"
$ws.Cells.Item(4, 3).Value = "Fzd4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.6836936666666666
$ws.Cells.Item(4, 8).Value = 2.051081
$ws.Cells.Item(4, 9).Value = 0.07240868516880868
$ws.Cells.Item(4, 10).Value = 0.07240868516880865
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 25.95900466666667
$ws.Cells.Item(4, 14).Value = 77.877014
$ws.Cells.Item(4, 15).Value = 0.4377032450206762
$ws.Cells.Item(4, 16).Value = 0.4377032450206762
$ws.Cells.Item(4, 17).Value = 17.74800708357045
$ws.Cells.Item(4, 18).Value = 159.732063752134
$ws.Cells.Item(4, 19).Value = 0.03169351646606806
$ws.Cells.Item(4, 20).Value = 0.03169351646606806

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Wnt9a"
$ws.Cells.Item(5, 3).Value = "Fzd4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7.8617
$ws.Cells.Item(5, 8).Value = 23.5851
$ws.Cells.Item(5, 9).Value = 0.8326175712099471
$ws.Cells.Item(5, 10).Value = 0.8326175712099468
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 13.57958433333333
$ws.Cells.Item(5, 14).Value = 40.738753
$ws.Cells.Item(5, 15).Value = 0.2289698008477291
$ws.Cells.Item(5, 16).Value = 0.2289698008477291
$ws.Cells.Item(5, 17).Value = 106.7586181533667
$ws.Cells.Item(5, 18).Value = 960.8275633803
$ws.Cells.Item(5, 19).Value = 0.1906442794622615
$ws.Cells.Item(5, 20).Value = 0.1906442794622614

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt9a"
$ws.Cells.Item(6, 3).Value = "Fzd4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.8617
$ws.Cells.Item(6, 8).Value = 23.5851
$ws.Cells.Item(6, 9).Value = 0.8326175712099471
$ws.Cells.Item(6, 10).Value = 0.8326175712099468
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 19.768727
$ws.Cells.Item(6, 14).Value = 59.306181
$ws.Cells.Item(6, 15).Value = 0.3333269541315948
$ws.Cells.Item(6, 16).Value = 0.3333269541315948
$ws.Cells.Item(6, 17).Value = 155.4158010559
$ws.Cells.Item(6, 18).Value = 1398.7422095031
$ws.Cells.Item(6, 19).Value = 0.2775338789678579
$ws.Cells.Item(6, 20).Value = 0.2775338789678578

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt9a"
$ws.Cells.Item(7, 3).Value = "Fzd4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.8617
$ws.Cells.Item(7, 8).Value = 23.5851
$ws.Cells.Item(7, 9).Value = 0.8326175712099471
$ws.Cells.Item(7, 10).Value = 0.8326175712099468
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 25.95900466666667
$ws.Cells.Item(7, 14).Value = 77.877014
$ws.Cells.Item(7, 15).Value = 0.4377032450206762
$ws.Cells.Item(7, 16).Value = 0.4377032450206762
$ws.Cells.Item(7, 17).Value = 204.0819069879333
$ws.Cells.Item(7, 18).Value = 1836.7371628914
$ws.Cells.Item(7, 19).Value = 0.3644394127798278
$ws.Cells.Item(7, 20).Value = 0.3644394127798277

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Wnt9a"
$ws.Cells.Item(8, 3).Value = "Fzd4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8967563333333333
$ws.Cells.Item(8, 8).Value = 2.690269
$ws.Cells.Item(8, 9).Value = 0.09497374362124447
$ws.Cells.Item(8, 10).Value = 0.09497374362124444
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 13.57958433333333
$ws.Cells.Item(8, 14).Value = 40.738753
$ws.Cells.Item(8, 15).Value = 0.2289698008477291
$ws.Cells.Item(8, 16).Value = 0.2289698008477291
$ws.Cells.Item(8, 17).Value = 12.17757825495078
$ws.Cells.Item(8, 18).Value = 109.598204294557
$ws.Cells.Item(8, 19).Value = 0.02174611916271962
$ws.Cells.Item(8, 20).Value = 0.02174611916271962

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Wnt9a"
$ws.Cells.Item(9, 3).Value = "Fzd4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8967563333333333
$ws.Cells.Item(9, 8).Value = 2.690269
$ws.Cells.Item(9, 9).Value = 0.09497374362124447
$ws.Cells.Item(9, 10).Value = 0.09497374362124444
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 19.768727
$ws.Cells.Item(9, 14).Value = 59.306181
$ws.Cells.Item(9, 15).Value = 0.3333269541315948
$ws.Cells.Item(9, 16).Value = 0.3333269541315948
$ws.Cells.Item(9, 17).Value = 17.72773113918766
$ws.Cells.Item(9, 18).Value = 159.549580252689
$ws.Cells.Item(9, 19).Value = 0.0316573086837444
$ws.Cells.Item(9, 20).Value = 0.03165730868374439

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Wnt9a"
$ws.Cells.Item(10, 3).Value = "Fzd4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8967563333333333
$ws.Cells.Item(10, 8).Value = 2.690269
$ws.Cells.Item(10, 9).Value = 0.09497374362124447
$ws.Cells.Item(10, 10).Value = 0.09497374362124444
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 25.95900466666667
$ws.Cells.Item(10, 14).Value = 77.877014
$ws.Cells.Item(10, 15).Value = 0.4377032450206762
$ws.Cells.Item(10, 16).Value = 0.4377032450206762
$ws.Cells.Item(10, 17).Value = 23.27890184186289
$ws.Cells.Item(10, 18).Value = 209.510116576766
$ws.Cells.Item(10, 19).Value = 0.04157031577478045
$ws.Cells.Item(10, 20).Value = 0.04157031577478044
